$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B-G (Forecast error table - naive component forecaster fix)
$data = @{
    2  = @(0.2796718369229467, 0.7060998478836559, 1.307254641040476, 1.143352369587118, 1.120870316762238, 46)
    3  = @(0.06641828055913036, 0.6381699131279182, 1.188041105879228, 1.08997298401347, 1.100241061446679, 45)
    4  = @(0.2050740221088844, 0.6171955800063507, 1.178333655589288, 1.085510780964099, 1.078287242720542, 44)
    5  = @(0.1096554176903501, 0.6078660076562977, 1.177810379744086, 1.08526972672423, 1.092493881120005, 43)
    6  = @(0.1717744169884318, 0.6286604205121966, 1.216341350767932, 1.102878665478634, 1.102625091828489, 42)
    7  = @(0.1768538724432516, 0.6104058027863382, 1.140351552268331, 1.067872441946289, 1.066208758438943, 41)
    8  = @(0.1598281167258759, 0.6402851931566111, 1.259767004438407, 1.122393426761938, 1.125108303637986, 40)
    9  = @(0.1949165218041688, 0.6346525613013413, 1.269321839014417, 1.126641841498183, 1.124158704767178, 39)
    10 = @(0.1797154758037312, 0.6670396664056156, 1.335955379238441, 1.155835359918722, 1.157104887042038, 38)
    11 = @(0.201330222591291, 0.6632936024058411, 1.352772756718701, 1.163087596322264, 1.16133112362294, 37)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}

$wb.Save()
